$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10 down to 11, 9 down to 10, ... 4 down to 5 (copy full format+content),
# working from the bottom up so we don't overwrite data before it is copied.
# (Paste formats, then values, separately -- pasting straight into a brand-new row
# with a single xlPasteAll drops styling on empty cells.)
for ($r = 10; $r -ge 4; $r--) {
    $src = $ws.Range("A" + $r + ":F" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":F" + ($r + 1))
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = 0
}

# --- Row 3: clear C3's value and give it the red-filled bordered style ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").ClearContents() | Out-Null
$ws.Range("C3").Interior.ColorIndex = 3
$excel.CutCopyMode = 0

# --- Build the new row 4 (blank record, only F4 carries a label) ---
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null

$ws.Range("B3").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null

$ws.Range("C3").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null

$ws.Range("D3").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null

$ws.Range("E3").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null

$ws.Range("F3").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A4:F4").ClearContents() | Out-Null
$ws.Range("F4").Value = "emailKosong"

# --- Update the active selection to match the new state ---
$ws.Range("G5").Select() | Out-Null
